$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(0.28864,0.28864,0.27809,0.27809,0.27809,0.27809,0.27685,0.27809,0.27498,0.2905,0.28305,0.28305,0.28305,0.28305,0.27312,0.27498,0.27685,0.26567,0.2905,0.2725,0.28802,0.2725,0.2843,0.28367,0.2843,0.27312,0.27809,0.2843,0.27809,0.26754,0.27747,0.28181,0.2905,0.28988,0.28988,0.28988,0.29981,0.29671,0.28988,0.28926,0.27747,0.27312,0.27312,0.28181,0.27623,0.27747,0.27623,0.27623,0.26567,0.27064,0.26878,0.26567,0.25947,0.25885,0.26567,0.25947,0.2545,0.25885,0.27498,0.27126,0.26505,0.25947,0.26505,0.26071,0.26133,0.26133,0.27188,0.26816,0.28367,0.2874,0.28988,0.26691,0.28119,0.27126,0.26878,0.28181,0.28181,0.28367,0.30292,0.29361,0.29671,0.30912,0.31533,0.31844,0.30912,0.29671,0.29671,0.2905,0.2905,0.30292,0.29671,0.29671,0.2905,0.29361,0.28492,0.28243,0.27188,0.26816,0.2843,0.30602,0.30602,0.29671,0.28802,0.2843,0.27933,0.28926,0.2905,0.30292,0.30912,0.29671,0.28554,0.2905,0.31533,0.30292,0.30912,0.29361,0.31533,0.35878,0.39603,0.39603,0.38361,0.38982,0.38051,0.38051,0.37741,0.40534,0.3712,0.3712,0.41775,0.45189,0.44569,0.4581)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = [double]$values[$i]
}

# Copy the existing column-A style (from A127) down into the newly added rows (A128:A133)
$ws.Range("A127").Copy() | Out-Null
$ws.Range("A128:A133").PasteSpecial(-4122) | Out-Null
